# Auto-generated Excel COM-interop edit script
# Applies numeric updates (and a few cell adds/removals) to the
# per-sheet "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H-N)
# as captured by the authoritative diff between before/after OOXML.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value2 = 6666.0527
$ws.Range("I32").Value2 = 7746.857
$ws.Range("J32").Value2 = 6035.5835
$ws.Range("K32").Value2 = 7746.857
$ws.Range("L32").Value2 = 6035.5835
$ws.Range("M32").Value2 = -7420.857
$ws.Range("N32").Value2 = -6687.5835
$ws.Range("H59").Value2 = 2399.5
$ws.Range("I59").Value2 = 1000
$ws.Range("K59").Value2 = 3000
$ws.Range("M59").Value2 = -2443
$ws.Range("H62").Value2 = 3826
$ws.Range("I62").Value2 = 3621
$ws.Range("J62").Value2 = 4543.5
$ws.Range("K62").Value2 = 3621
$ws.Range("L62").Value2 = 4543.5
$ws.Range("M62").Value2 = -2997
$ws.Range("N62").Value2 = -5791.5
$ws.Range("H65").Value2 = 3826
$ws.Range("I65").Value2 = 3621
$ws.Range("J65").Value2 = 4543.5
$ws.Range("K65").Value2 = 18105
$ws.Range("L65").Value2 = 22717.5
$ws.Range("M65").Value2 = -14985
$ws.Range("N65").Value2 = -28957.5
$ws.Range("H74").Value2 = 19437.5
$ws.Range("I74").Value2 = 2750
$ws.Range("K74").Value2 = 2750
$ws.Range("M74").Value2 = -1814
$ws.Range("H76").Value2 = 2540.6
$ws.Range("I76").Value2 = 2149.5
$ws.Range("J76").Value2 = 2801.3333
$ws.Range("K76").Value2 = 2149.5
$ws.Range("L76").Value2 = 2801.3333
$ws.Range("M76").Value2 = -1834.5
$ws.Range("N76").Value2 = -3431.3333
$ws.Range("H77").Value2 = 19437.5
$ws.Range("I77").Value2 = 2750
$ws.Range("K77").Value2 = 13750
$ws.Range("M77").Value2 = -9070
$ws.Range("H79").Value2 = 2540.6
$ws.Range("I79").Value2 = 2149.5
$ws.Range("J79").Value2 = 2801.3333
$ws.Range("K79").Value2 = 2149.5
$ws.Range("L79").Value2 = 2801.3333
$ws.Range("M79").Value2 = -1057.5
$ws.Range("N79").Value2 = -4985.3333
$ws.Range("H112").Value2 = 1859.3572
$ws.Range("I112").Value2 = 1423.5
$ws.Range("J112").Value2 = 1932
$ws.Range("K112").Value2 = 4270.5
$ws.Range("L112").Value2 = 5796
$ws.Range("M112").Value2 = -3162.5
$ws.Range("N112").Value2 = -8012
$ws.Range("H137").Value2 = 970
$ws.Range("I137").Value2 = 970
$ws.Range("J137").Value2 = 0
$ws.Range("K137").Value2 = 2910
$ws.Range("L137").Value2 = 0
$ws.Range("M137").Value2 = -360
$ws.Range("N137").ClearContents()
$ws.Range("H141").Value2 = 4500.9443
$ws.Range("J141").Value2 = 3494.5
$ws.Range("L141").Value2 = 10483.5
$ws.Range("N141").Value2 = -20843.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value2 = 1680.1017
$ws.Range("I132").Value2 = 1660.804
$ws.Range("J132").Value2 = 1803.125
$ws.Range("K132").Value2 = 4982.412
$ws.Range("L132").Value2 = 5409.375
$ws.Range("M132").Value2 = -2452.412
$ws.Range("N132").Value2 = -10469.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 26925474
$ws.Range("I58").Value2 = 14287687
$ws.Range("K58").Value2 = 14287687
$ws.Range("M58").Value2 = -14287484
$ws.Range("H122").Value2 = 981.7273
$ws.Range("I122").Value2 = 985.75
$ws.Range("J122").Value2 = 971
$ws.Range("K122").Value2 = 2957.25
$ws.Range("L122").Value2 = 2913
$ws.Range("M122").Value2 = -507.25
$ws.Range("N122").Value2 = -7813
$ws.Range("H132").Value2 = 1743.2413
$ws.Range("I132").Value2 = 1748.3572
$ws.Range("K132").Value2 = 5245.071599999999
$ws.Range("M132").Value2 = -2715.071599999999
$ws.Range("H134").Value2 = 2166.8333
$ws.Range("I134").Value2 = 2166.8333
$ws.Range("K134").Value2 = 6500.499899999999
$ws.Range("M134").Value2 = -3965.499899999999
$ws.Range("H136").Value2 = 26925474
$ws.Range("I136").Value2 = 14287687
$ws.Range("K136").Value2 = 42863061
$ws.Range("M136").Value2 = -42860511

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 2197.8125
$ws.Range("I5").Value2 = 839.5
$ws.Range("K5").Value2 = 2518.5
$ws.Range("M5").Value2 = -2406.5
$ws.Range("H7").Value2 = 682.125
$ws.Range("I7").Value2 = 719.5714
$ws.Range("K7").Value2 = 2158.7142
$ws.Range("M7").Value2 = -2046.7142
$ws.Range("H56").Value2 = 13933.267
$ws.Range("I56").Value2 = 13933.267
$ws.Range("K56").Value2 = 13933.267
$ws.Range("M56").Value2 = -13403.267
$ws.Range("H63").Value2 = 2013.1765
$ws.Range("J63").Value2 = 1992.8334
$ws.Range("L63").Value2 = 5978.5002
$ws.Range("N63").Value2 = -7476.5002
$ws.Range("H64").Value2 = 6584.375
$ws.Range("I64").Value2 = 3637
$ws.Range("J64").Value2 = 8352.799999999999
$ws.Range("K64").Value2 = 10911
$ws.Range("L64").Value2 = 25058.4
$ws.Range("M64").Value2 = -10641
$ws.Range("N64").Value2 = -25598.4
$ws.Range("H66").Value2 = 2013.1765
$ws.Range("J66").Value2 = 1992.8334
$ws.Range("L66").Value2 = 17935.5006
$ws.Range("N66").Value2 = -25423.5006
$ws.Range("H67").Value2 = 6584.375
$ws.Range("I67").Value2 = 3637
$ws.Range("J67").Value2 = 8352.799999999999
$ws.Range("K67").Value2 = 10911
$ws.Range("L67").Value2 = 25058.4
$ws.Range("M67").Value2 = -9975
$ws.Range("N67").Value2 = -26930.4
$ws.Range("H70").Value2 = 1427.5
$ws.Range("J70").Value2 = 0
$ws.Range("L70").Value2 = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value2 = 1427.5
$ws.Range("J73").Value2 = 0
$ws.Range("L73").Value2 = 0
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value2 = 9000
$ws.Range("J74").Value2 = 9000
$ws.Range("L74").Value2 = 27000
$ws.Range("N74").Value2 = -29122
$ws.Range("H75").Value2 = 2332.2354
$ws.Range("J75").Value2 = 2330.8
$ws.Range("L75").Value2 = 6992.400000000001
$ws.Range("N75").Value2 = -8988.400000000001
$ws.Range("H76").Value2 = 2999.5
$ws.Range("I76").Value2 = 2999.5
$ws.Range("J76").Value2 = 0
$ws.Range("K76").Value2 = 8998.5
$ws.Range("L76").Value2 = 0
$ws.Range("M76").Value2 = -8615.5
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value2 = 9000
$ws.Range("J77").Value2 = 9000
$ws.Range("L77").Value2 = 81000
$ws.Range("N77").Value2 = -91608
$ws.Range("H78").Value2 = 2332.2354
$ws.Range("J78").Value2 = 2330.8
$ws.Range("L78").Value2 = 20977.2
$ws.Range("N78").Value2 = -30961.2
$ws.Range("H79").Value2 = 2999.5
$ws.Range("I79").Value2 = 2999.5
$ws.Range("J79").Value2 = 0
$ws.Range("K79").Value2 = 8998.5
$ws.Range("L79").Value2 = 0
$ws.Range("M79").Value2 = -7672.5
$ws.Range("N79").ClearContents()
$ws.Range("H82").Value2 = 7666
$ws.Range("J82").Value2 = 8499
$ws.Range("L82").Value2 = 25497
$ws.Range("N82").Value2 = -26309
$ws.Range("H85").Value2 = 7666
$ws.Range("J85").Value2 = 8499
$ws.Range("L85").Value2 = 25497
$ws.Range("N85").Value2 = -28305
$ws.Range("H87").Value2 = 6610.4
$ws.Range("I87").Value2 = 6388.25
$ws.Range("J87").Value2 = 7499
$ws.Range("K87").Value2 = 19164.75
$ws.Range("L87").Value2 = 22497
$ws.Range("M87").Value2 = -17916.75
$ws.Range("N87").Value2 = -24993
$ws.Range("H88").Value2 = 3000
$ws.Range("J88").Value2 = 3000
$ws.Range("L88").Value2 = 9000
$ws.Range("N88").Value2 = -9856
$ws.Range("H90").Value2 = 6610.4
$ws.Range("I90").Value2 = 6388.25
$ws.Range("J90").Value2 = 7499
$ws.Range("K90").Value2 = 57494.25
$ws.Range("L90").Value2 = 67491
$ws.Range("M90").Value2 = -51254.25
$ws.Range("N90").Value2 = -79971
$ws.Range("H91").Value2 = 3000
$ws.Range("J91").Value2 = 3000
$ws.Range("L91").Value2 = 9000
$ws.Range("N91").Value2 = -11964
$ws.Range("H92").Value2 = 498.875
$ws.Range("J92").Value2 = 541.5714
$ws.Range("L92").Value2 = 1624.7142
$ws.Range("N92").Value2 = -4120.7142
$ws.Range("H124").Value2 = 709.6667
$ws.Range("I124").Value2 = 709.6667
$ws.Range("K124").Value2 = 2129.0001
$ws.Range("M124").Value2 = 2780.9999
$ws.Range("H129").Value2 = 6093.273
$ws.Range("I129").Value2 = 1328.375
$ws.Range("J129").Value2 = 18799.666
$ws.Range("K129").Value2 = 3985.125
$ws.Range("L129").Value2 = 56398.99800000001
$ws.Range("M129").Value2 = 1014.875
$ws.Range("N129").Value2 = -66398.99800000001
$ws.Range("H135").Value2 = 2197.8125
$ws.Range("I135").Value2 = 839.5
$ws.Range("K135").Value2 = 7555.5
$ws.Range("M135").Value2 = -5020.5
$ws.Range("H139").Value2 = 1902.125
$ws.Range("I139").Value2 = 1461.3334
$ws.Range("K139").Value2 = 4384.0002
$ws.Range("M139").Value2 = 755.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 76436.516
$ws.Range("I122").Value2 = 74345.03999999999
$ws.Range("K122").Value2 = 223035.12
$ws.Range("M122").Value2 = -220585.12
$ws.Range("H132").Value2 = 5633.567
$ws.Range("I132").Value2 = 5036
$ws.Range("J132").Value2 = 13999.5
$ws.Range("K132").Value2 = 15108
$ws.Range("L132").Value2 = 41998.5
$ws.Range("M132").Value2 = -12578
$ws.Range("N132").Value2 = -47058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 1137.875
$ws.Range("I16").Value2 = 1254.5385
$ws.Range("J16").Value2 = 632.3333
$ws.Range("K16").Value2 = 1254.5385
$ws.Range("L16").Value2 = 632.3333
$ws.Range("M16").Value2 = -1084.5385
$ws.Range("N16").Value2 = -972.3333
$ws.Range("H132").Value2 = 3799.9
$ws.Range("I132").Value2 = 3416.5557
$ws.Range("J132").Value2 = 7250
$ws.Range("K132").Value2 = 10249.6671
$ws.Range("L132").Value2 = 21750
$ws.Range("M132").Value2 = -7719.667099999999
$ws.Range("N132").Value2 = -26810

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value2 = 0
$ws.Range("J94").Value2 = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").ClearContents()
